$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.146.71"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +3.15%  '
$ws.Range("D3").Value = "'3.263.38"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +2.62%  '
$ws.Range("D5").Value = "'546.76"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.71%  '
$ws.Range("D6").Value = "'148.47"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +4.80%  '
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.29%  '
$ws.Range("D8").Value = "'0.524"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.73%  '
$ws.Range("D9").Value = "'7.42"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.92%  '
$ws.Range("E10").Value = '  +3.89%  '
$ws.Range("D11").Value = "'0.434"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.29%  '
$ws.Range("D12").Value = "'3.822.71"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +2.42%  '
$ws.Range("E13").Value = '  -1.19%  '
$ws.Range("D14").Value = "'26.64"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +2.95%  '
$ws.Range("E15").Value = '  +3.58%  '
$ws.Range("D16").Value = "'61.117.50"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.99%  '
$ws.Range("D17").Value = "'3.264.51"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +2.67%  '
$ws.Range("D18").Value = "'6.33"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.51%  '
$ws.Range("D19").Value = "'13.51"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +3.95%  '
$ws.Range("E20").Value = '  +3.66%  '
$ws.Range("D21").Value = "'379.34"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.96%  '
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("D23").Value = "'0.535"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.50%  '
$ws.Range("D24").Value = "'70.21"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.54%  '
$ws.Range("E25").Value = '  +1.95%  '
$ws.Range("E26").Value = '  +2.91%  '
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.03%  '
$ws.Range("E28").Value = '  +6.63%  '
$ws.Range("E29").Value = '  +2.84%  '
$ws.Range("D30").Value = "'22.64"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.00%  '
$ws.Range("D31").Value = "'6.24"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +2.84%  '
$ws.Range("D32").Value = "'5.44"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.40%  '
$ws.Range("D33").Value = "'1.27"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +8.13%  '
$ws.Range("E34").Value = '  +4.94%  '
$ws.Range("D35").Value = "'159.61"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +1.56%  '
$ws.Range("E36").Value = '  +7.99%  '
$ws.Range("D37").Value = "'26.49"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +4.04%  '
$ws.Range("D38").Value = "'2.811.87"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +3.84%  '
$ws.Range("E39").Value = '  +1.12%  '
$ws.Range("E40").Value = '  +1.93%  '
$ws.Range("E41").Value = '  +6.92%  '
$ws.Range("E42").Value = '  +0.63%  '
$ws.Range("D43").Value = "'40.20"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.62%  '
$ws.Range("D44").Value = "'0.735"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.17%  '
$ws.Range("D45").Value = "'3.305.16"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +2.54%  '
$ws.Range("E46").Value = '  +2.58%  '
$ws.Range("E47").Value = '  +2.85%  '
$ws.Range("D48").Value = "'21.64"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +6.40%  '
$ws.Range("E49").Value = '  +0.78%  '
$ws.Range("D50").Value = "'0.807"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +5.97%  '
$ws.Range("D51").Value = "'278.93"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +8.76%  '
